$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/25/2025  Through  8/31/2025"

# --- Cells changing type (number <-> text placeholder) ---
# Strategy: Copy format+value from a stable donor cell that already carries the exact
# target style (this engine creates a brand-new style/numFmt entry whenever NumberFormat
# is set directly via COM, so reusing an existing donor cell keeps the style table intact).
# Donors (unmodified elsewhere in this script):
#   C14 = style 13, shared text "0"      E14 = style 13, shared text "***.*"
#   I14 = style 14, plain number          K14 = style 15, plain number (% fmt)
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("I14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 3
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("I14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("I14").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K14").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100

# --- Pure value updates (style/type unchanged) ---
$ws.Range("N14").Value = -75
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("M15").Value = -4.545454545454
$ws.Range("N15").Value = -70.422535211267
$ws.Range("C16").Value = 6
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 21
$ws.Range("H16").Value = 31.25
$ws.Range("I16").Value = 186
$ws.Range("J16").Value = 182
$ws.Range("K16").Value = 2.197802197802
$ws.Range("L16").Value = 12.048192771084
$ws.Range("M16").Value = -28.735632183908
$ws.Range("N16").Value = -88.411214953271
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = -58.333333333333
$ws.Range("F17").Value = 58
$ws.Range("G17").Value = 88
$ws.Range("H17").Value = -34.090909090909
$ws.Range("I17").Value = 501
$ws.Range("J17").Value = 571
$ws.Range("K17").Value = -12.259194395796
$ws.Range("L17").Value = 8.676789587852
$ws.Range("M17").Value = 65.346534653465
$ws.Range("N17").Value = -42.677345537757
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 6.666666666666
$ws.Range("I18").Value = 121
$ws.Range("J18").Value = 139
$ws.Range("K18").Value = -12.949640287769
$ws.Range("L18").Value = -7.633587786259
$ws.Range("M18").Value = -51.405622489959
$ws.Range("N18").Value = -91.050295857988
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 112.5
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 22.916666666666
$ws.Range("I19").Value = 457
$ws.Range("J19").Value = 399
$ws.Range("K19").Value = 14.536340852130
$ws.Range("L19").Value = 3.863636363636
$ws.Range("M19").Value = 28.370786516853
$ws.Range("N19").Value = -9.504950495049
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 18.75
$ws.Range("I20").Value = 134
$ws.Range("J20").Value = 137
$ws.Range("K20").Value = -2.189781021897
$ws.Range("L20").Value = -14.102564102564
$ws.Range("M20").Value = -30.208333333333
$ws.Range("N20").Value = -90.037174721189
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = -4.444444444444
$ws.Range("F21").Value = 175
$ws.Range("G21").Value = 185
$ws.Range("H21").Value = -5.405405405405
$ws.Range("I21").Value = 1427
$ws.Range("J21").Value = 1461
$ws.Range("K21").Value = -2.327173169062
$ws.Range("L21").Value = 2.661870503597
$ws.Range("M21").Value = 2.074391988555
$ws.Range("N21").Value = -75.311418685121
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -33.333333333333
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -4.347826086956
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 28.571428571428
$ws.Range("I24").Value = 742
$ws.Range("J24").Value = 827
$ws.Range("K24").Value = -10.278113663845
$ws.Range("L24").Value = -4.993597951344
$ws.Range("M24").Value = 16.300940438871
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 108.333333333333
$ws.Range("I25").Value = 133
$ws.Range("J25").Value = 132
$ws.Range("K25").Value = 0.757575757575
$ws.Range("L25").Value = -5.673758865248
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -6.25
$ws.Range("F26").Value = 53
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = -26.388888888888
$ws.Range("I26").Value = 608
$ws.Range("J26").Value = 633
$ws.Range("K26").Value = -3.949447077409
$ws.Range("L26").Value = 19.685039370078
$ws.Range("M26").Value = -3.492063492063
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = -35.135135135135
$ws.Range("L27").Value = -36.842105263157
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 58
$ws.Range("K28").Value = 23.404255319148
$ws.Range("L28").Value = 23.404255319148
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 250
$ws.Range("J29").Value = 25
$ws.Range("K29").Value = 28
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = -78.947368421052
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 100
$ws.Range("J30").Value = 21
$ws.Range("K30").Value = 14.285714285714
$ws.Range("L30").Value = -11.111111111111
$ws.Range("N30").Value = -82.733812949640
